# Auto-generated Excel COM-interop edit script
# Updates crypto price/volume figures and re-orders coin rows per upstream diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '''91.958.52'
$ws.Range("E2").Value2 = '''  -6.16%  '
$ws.Range("D3").Value2 = '''3.306.68'
$ws.Range("E3").Value2 = '''  -5.28%  '
$ws.Range("D4").Value2 = '''1.00'
$ws.Range("E4").Value2 = '''  +0.05%  '
$ws.Range("D5").Value2 = '''226.61'
$ws.Range("E5").Value2 = '''  -9.79%  '
$ws.Range("D6").Value2 = '''617.94'
$ws.Range("E6").Value2 = '''  -7.19%  '
$ws.Range("D7").Value2 = '''1.32'
$ws.Range("E7").Value2 = '''  -10.30%  '
$ws.Range("D8").Value2 = '''0.370'
$ws.Range("E8").Value2 = '''  -12.79%  '
$ws.Range("D9").Value2 = '''1.00'
$ws.Range("E9").Value2 = '''  +0.15%  '
$ws.Range("D10").Value2 = '''0.892'
$ws.Range("E10").Value2 = '''  -14.16%  '
$ws.Range("D11").Value2 = '''3.304.56'
$ws.Range("E11").Value2 = '''  -5.28%  '
$ws.Range("D12").Value2 = '''41.13'
$ws.Range("E12").Value2 = '''  -8.85%  '
$ws.Range("D13").Value2 = '''0.188'
$ws.Range("E13").Value2 = '''  -10.59%  '
$ws.Range("D14").Value2 = '''91.790.92'
$ws.Range("E14").Value2 = '''  -6.19%  '
$ws.Range("D15").Value2 = '''5.84'
$ws.Range("E15").Value2 = '''  -6.14%  '
$ws.Range("D16").Value2 = '''3.931.14'
$ws.Range("E16").Value2 = '''  -5.31%  '
$ws.Range("D17").Value2 = '''0.0000236'
$ws.Range("E17").Value2 = '''  -9.07%  '
$ws.Range("D18").Value2 = '''7.78'
$ws.Range("E18").Value2 = '''  -12.44%  '
$ws.Range("D19").Value2 = '''3.304.84'
$ws.Range("E19").Value2 = '''  -5.32%  '
$ws.Range("D20").Value2 = '''16.52'
$ws.Range("E20").Value2 = '''  -12.04%  '
$ws.Range("D21").Value2 = '''10.68'
$ws.Range("E21").Value2 = '''  -10.70%  '
$ws.Range("D22").Value2 = '''480.18'
$ws.Range("E22").Value2 = '''  -8.02%  '
$ws.Range("D23").Value2 = '''3.19'
$ws.Range("E23").Value2 = '''  -5.88%  '
$ws.Range("D24").Value2 = '''0.429'
$ws.Range("E24").Value2 = '''  -16.52%  '
$ws.Range("D25").Value2 = '''0.0000176'
$ws.Range("E25").Value2 = '''  -12.30%  '
$ws.Range("D26").Value2 = '''5.97'
$ws.Range("E26").Value2 = '''  -12.34%  '
$ws.Range("D27").Value2 = '''88.55'
$ws.Range("E27").Value2 = '''  -9.17%  '
$ws.Range("B28").Value2 = 'WrappedeETH'
$ws.Range("C28").Value2 = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").Value2 = '''3.494.49'
$ws.Range("E28").Value2 = '''  -5.09%  '
$ws.Range("B29").Value2 = 'Aptos'
$ws.Range("C29").Value2 = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D29").Value2 = '''11.24'
$ws.Range("E29").Value2 = '''  -11.51%  '
$ws.Range("B30").Value2 = 'Dai'
$ws.Range("C30").Value2 = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D30").Value2 = '''1.00'
$ws.Range("E30").Value2 = '''  +0.06%  '
$ws.Range("B31").Value2 = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value2 = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value2 = '''10.86'
$ws.Range("E31").Value2 = '''  -12.42%  '
$ws.Range("D32").Value2 = '''0.132'
$ws.Range("E32").Value2 = '''  -7.96%  '
$ws.Range("B33").Value2 = 'PancakeSwap'
$ws.Range("C33").Value2 = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D33").Value2 = '''2.55'
$ws.Range("E33").Value2 = '''  -11.39%  '
$ws.Range("B34").Value2 = 'Binance-PegBSC-USD'
$ws.Range("C34").Value2 = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D34").Value2 = '''0.997'
$ws.Range("E34").Value2 = '''  -0.32%  '
$ws.Range("B35").Value2 = 'Cronos'
$ws.Range("C35").Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D35").Value2 = '''0.166'
$ws.Range("E35").Value2 = '''  -12.24%  '
$ws.Range("B36").Value2 = 'EthereumClassic'
$ws.Range("C36").Value2 = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value2 = '''27.76'
$ws.Range("E36").Value2 = '''  -11.00%  '
$ws.Range("B37").Value2 = 'PolygonEcosystemToken'
$ws.Range("C37").Value2 = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D37").Value2 = '''0.510'
$ws.Range("E37").Value2 = '''  -14.44%  '
$ws.Range("D38").Value2 = '''520.97'
$ws.Range("E38").Value2 = '''  -0.62%  '
$ws.Range("B39").Value2 = 'USDe'
$ws.Range("C39").Value2 = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D39").Value2 = '''1.00'
$ws.Range("E39").Value2 = '''  -0.03%  '
$ws.Range("B40").Value2 = 'RenderToken'
$ws.Range("C40").Value2 = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").Value2 = '''7.16'
$ws.Range("E40").Value2 = '''  -8.95%  '
$ws.Range("E41").Value2 = '''  -8.43%  '
$ws.Range("B42").Value2 = 'Fetch.AI'
$ws.Range("C42").Value2 = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D42").Value2 = '''1.33'
$ws.Range("E42").Value2 = '''  -11.69%  '
$ws.Range("B43").Value2 = 'ARBITRUM'
$ws.Range("C43").Value2 = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D43").Value2 = '''0.855'
$ws.Range("E43").Value2 = '''  -6.80%  '
$ws.Range("B44").Value2 = 'WhiteBITCoin'
$ws.Range("C44").Value2 = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D44").Value2 = '''23.91'
$ws.Range("E44").Value2 = '''  -1.91%  '
$ws.Range("B45").Value2 = 'MantraDAO'
$ws.Range("C45").Value2 = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D45").Value2 = '''3.52'
$ws.Range("E45").Value2 = '''  -3.96%  '
$ws.Range("B46").Value2 = 'ImmutableX'
$ws.Range("C46").Value2 = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D46").Value2 = '''1.62'
$ws.Range("E46").Value2 = '''  -6.41%  '
$ws.Range("B47").Value2 = 'Filecoin'
$ws.Range("C47").Value2 = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D47").Value2 = '''5.21'
$ws.Range("E47").Value2 = '''  -9.14%  '
$ws.Range("B48").Value2 = 'Stacks'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value2 = '''2.07'
$ws.Range("E48").Value2 = '''  -7.37%  '
$ws.Range("B49").Value2 = 'OKB'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value2 = '''51.78'
$ws.Range("E49").Value2 = '''  -5.23%  '
$ws.Range("B50").Value2 = 'VeChain'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value2 = '''0.0381'
$ws.Range("E50").Value2 = '''  -12.06%  '
$ws.Range("B51").Value2 = 'Cosmos'
$ws.Range("C51").Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value2 = '''7.74'
$ws.Range("E51").Value2 = '''  -11.35%  '
